# Edit slide 3 ("Notable Technologies Used") content placeholder:
#  - switch bodyPr autofit from normAutofit(lnSpcReduction=10000) to plain normAutofit
#  - remove the "RabbitMQ" paragraph
#  - remove the stray empty paragraph that followed "Docker"
#  - shift the paragraph-targeted animations (Trello/Slack/GitHub) back by 2
#    to account for the two removed paragraphs

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# 1. Reduce the autofit line-spacing reduction back to a plain normAutofit.
$tf.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

# 2. Remove the "RabbitMQ" paragraph entirely (currently paragraph 4).
$tr.Paragraphs(4, 1).Delete()

# 3. Merge away the now-adjacent empty paragraph that trailed "Docker"
#    (currently paragraphs 4 "Docker" + 5 "") without leaving a stray
#    endParaRPr behind on the "Docker" paragraph.
$tr.Paragraphs(4, 2).Text = "Docker"

# 4. Re-target the three click-triggered "Fly In" animations on this shape
#    so they point at the same visible bullets (Trello/Slack/GitHub) now
#    that two paragraphs were removed ahead of them.
$main = $s.TimeLine.MainSequence
$main.Item(1).Paragraph = 5
$main.Item(2).Paragraph = 6
$main.Item(3).Paragraph = 7
